$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"0"
$ws.Cells.Item(2, 2).Value = [double]"0"
$ws.Cells.Item(2, 3).Value = [double]"2.5"
$ws.Cells.Item(2, 4).Value = [double]"-0.5824953833978498"
$ws.Cells.Item(2, 5).Value = [double]"-33.37452705455156"
$ws.Cells.Item(2, 6).Value = [double]"0"
$ws.Cells.Item(2, 7).Value = [double]"-3.109580457508859e-40"
$ws.Cells.Item(2, 8).Value = [double]"0"

$ws.Cells.Item(3, 1).Value = [double]"0.04899084750364537"
$ws.Cells.Item(3, 2).Value = [double]"0.1725862640140226"
$ws.Cells.Item(3, 3).Value = [double]"2.386310384199221"
$ws.Cells.Item(3, 4).Value = [double]"-0.5824953833978498"
$ws.Cells.Item(3, 5).Value = [double]"-33.37452705455156"
$ws.Cells.Item(3, 6).Value = [double]"-1.853886236141795e-25"
$ws.Cells.Item(3, 7).Value = [double]"-1.704228138430471e-33"
$ws.Cells.Item(3, 8).Value = [double]"-2.353646874159769e-08"

$ws.Cells.Item(4, 1).Value = [double]"0.09798169500729074"
$ws.Cells.Item(4, 2).Value = [double]"0.3451725280280452"
$ws.Cells.Item(4, 3).Value = [double]"2.272620768398441"
$ws.Cells.Item(4, 4).Value = [double]"-0.5824953833978498"
$ws.Cells.Item(4, 5).Value = [double]"-33.37452705455156"
$ws.Cells.Item(4, 6).Value = [double]"-9.056617722904067e-09"
$ws.Cells.Item(4, 7).Value = [double]"4.925721685621824e-22"
$ws.Cells.Item(4, 8).Value = [double]"-3.537075797554134e-08"

$ws.Cells.Item(5, 1).Value = [double]"0.1469725425109361"
$ws.Cells.Item(5, 2).Value = [double]"0.5177587920420679"
$ws.Cells.Item(5, 3).Value = [double]"2.158931152597662"
$ws.Cells.Item(5, 4).Value = [double]"-0.5824953833978497"
$ws.Cells.Item(5, 5).Value = [double]"-33.37452705455155"
$ws.Cells.Item(5, 6).Value = [double]"-8.989233274530135e-09"
$ws.Cells.Item(5, 7).Value = [double]"-2.542164335276259e-08"
$ws.Cells.Item(5, 8).Value = [double]"-3.502387478029477e-08"

$ws.Cells.Item(6, 1).Value = [double]"0.1959633900145815"
$ws.Cells.Item(6, 2).Value = [double]"0.6903450560560904"
$ws.Cells.Item(6, 3).Value = [double]"2.045241536796882"
$ws.Cells.Item(6, 4).Value = [double]"-0.5824953834483982"
$ws.Cells.Item(6, 5).Value = [double]"-33.37452705744776"
$ws.Cells.Item(6, 6).Value = [double]"-8.890138497498464e-09"
$ws.Cells.Item(6, 7).Value = [double]"-5.065414066160378e-08"
$ws.Cells.Item(6, 8).Value = [double]"-3.450789491338028e-08"

$ws.Cells.Item(7, 1).Value = [double]"0.2449542375182269"
$ws.Cells.Item(7, 2).Value = [double]"0.8629313200643662"
$ws.Cells.Item(7, 3).Value = [double]"1.931551920987378"
$ws.Cells.Item(7, 4).Value = [double]"-0.5824953835491191"
$ws.Cells.Item(7, 5).Value = [double]"-33.37452706321865"
$ws.Cells.Item(7, 6).Value = [double]"-8.742156963774413e-09"
$ws.Cells.Item(7, 7).Value = [double]"-7.560848202347024e-08"
$ws.Cells.Item(7, 8).Value = [double]"-3.370950004657494e-08"

$ws.Cells.Item(8, 1).Value = [double]"0.2939450850218722"
$ws.Cells.Item(8, 2).Value = [double]"1.035517584061191"
$ws.Cells.Item(8, 3).Value = [double]"1.817862305160491"
$ws.Cells.Item(8, 4).Value = [double]"-0.5824953836994594"
$ws.Cells.Item(8, 5).Value = [double]"-33.37452707183251"
$ws.Cells.Item(8, 6).Value = [double]"-8.51044903591456e-09"
$ws.Cells.Item(8, 7).Value = [double]"-1.001474438377871e-07"
$ws.Cells.Item(8, 8).Value = [double]"-3.23699623947412e-08"

$ws.Cells.Item(9, 1).Value = [double]"0.3429359325255176"
$ws.Cells.Item(9, 2).Value = [double]"1.208103848040924"
$ws.Cells.Item(9, 3).Value = [double]"1.704172689307658"
$ws.Cells.Item(9, 4).Value = [double]"-0.582495383898593"
$ws.Cells.Item(9, 5).Value = [double]"-33.37452708324203"
$ws.Cells.Item(9, 6).Value = [double]"-8.113235445122279e-09"
$ws.Cells.Item(9, 7).Value = [double]"-1.240360087288136e-07"
$ws.Cells.Item(9, 8).Value = [double]"-2.974180643070676e-08"

$ws.Cells.Item(10, 1).Value = [double]"0.391926780029163"
$ws.Cells.Item(10, 2).Value = [double]"1.380690111998017"
$ws.Cells.Item(10, 3).Value = [double]"1.590483073420457"
$ws.Cells.Item(10, 4).Value = [double]"-0.5824953841452267"
$ws.Cells.Item(10, 5).Value = [double]"-33.3745270973731"
$ws.Cells.Item(10, 6).Value = [double]"-7.304622063497112e-09"
$ws.Cells.Item(10, 7).Value = [double]"-1.468096074651355e-07"
$ws.Cells.Item(10, 8).Value = [double]"-2.230156706816783e-08"

$ws.Cells.Item(11, 1).Value = [double]"0.4409176275328083"
$ws.Cells.Item(11, 2).Value = [double]"1.553276375927071"
$ws.Cells.Item(11, 3).Value = [double]"1.47679345749069"
$ws.Cells.Item(11, 4).Value = [double]"-0.5824953844371435"
$ws.Cells.Item(11, 5).Value = [double]"-33.3745271140987"
$ws.Cells.Item(11, 6).Value = [double]"-4.85427847843017e-09"
$ws.Cells.Item(11, 7).Value = [double]"-1.673134536711984e-07"
$ws.Cells.Item(11, 8).Value = [double]"1.949753791290979"

$ws.Cells.Item(12, 1).Value = [double]"0.4899084750364537"
$ws.Cells.Item(12, 2).Value = [double]"1.725862639822937"
$ws.Cells.Item(12, 3).Value = [double]"1.363103841510543"
$ws.Cells.Item(12, 4).Value = [double]"-0.5824953847698302"
$ws.Cells.Item(12, 5).Value = [double]"-33.37452713316024"
$ws.Cells.Item(12, 6).Value = [double]"0.7502474136425195"
$ws.Cells.Item(12, 7).Value = [double]"-1.809392618944799e-07"
$ws.Cells.Item(12, 8).Value = [double]"4.893077378402382"

$ws.Cells.Item(13, 1).Value = [double]"0.5388993225400991"
$ws.Cells.Item(13, 2).Value = [double]"1.898448903680979"
$ws.Cells.Item(13, 3).Value = [double]"1.249414225472978"
$ws.Cells.Item(13, 4).Value = [double]"-0.5824953851296106"
$ws.Cells.Item(13, 5).Value = [double]"-33.37452715377415"
$ws.Cells.Item(13, 6).Value = [double]"1.500000010289458"
$ws.Cells.Item(13, 7).Value = [double]"2.105920898981103"
$ws.Cells.Item(13, 8).Value = [double]"5.887276959717616"

$ws.Cells.Item(14, 1).Value = [double]"0.5878901700437444"
$ws.Cells.Item(14, 2).Value = [double]"2.071035167498119"
$ws.Cells.Item(14, 3).Value = [double]"1.13572460937332"
$ws.Cells.Item(14, 4).Value = [double]"-0.578306076748911"
$ws.Cells.Item(14, 5).Value = [double]"-33.13449746448127"
$ws.Cells.Item(14, 6).Value = [double]"1.50000001275109"
$ws.Cells.Item(14, 7).Value = [double]"6.316374122955148"
$ws.Cells.Item(14, 8).Value = [double]"5.887276965113485"

$ws.Cells.Item(15, 1).Value = [double]"0.6368810175473898"
$ws.Cells.Item(15, 2).Value = [double]"2.244096196315094"
$ws.Cells.Item(15, 3).Value = [double]"1.022759005882449"
$ws.Cells.Item(15, 4).Value = [double]"-0.5656954447792922"
$ws.Cells.Item(15, 5).Value = [double]"-32.41196147562936"
$ws.Cells.Item(15, 6).Value = [double]"1.500000013571332"
$ws.Cells.Item(15, 7).Value = [double]"10.52682735383892"
$ws.Cells.Item(15, 8).Value = [double]"5.88727696726643"

$ws.Cells.Item(16, 1).Value = [double]"0.6858718650510353"
$ws.Cells.Item(16, 2).Value = [double]"2.418567994429889"
$ws.Cells.Item(16, 3).Value = [double]"0.9119847357208399"
$ws.Cells.Item(16, 4).Value = [double]"-0.5445251126856346"
$ws.Cells.Item(16, 5).Value = [double]"-31.19899079577243"
$ws.Cells.Item(16, 6).Value = [double]"1.50000001398124"
$ws.Cells.Item(16, 7).Value = [double]"14.73728058702508"
$ws.Cells.Item(16, 8).Value = [double]"5.887276968448739"

$ws.Cells.Item(17, 1).Value = [double]"0.7348627125546806"
$ws.Cells.Item(17, 2).Value = [double]"2.595345649277633"
$ws.Cells.Item(17, 3).Value = [double]"0.8049286382122038"
$ws.Cells.Item(17, 4).Value = [double]"-0.5145576370584923"
$ws.Cells.Item(17, 5).Value = [double]"-29.48198091967602"
$ws.Cells.Item(17, 6).Value = [double]"1.500000014227027"
$ws.Cells.Item(17, 7).Value = [double]"18.94773382136184"
$ws.Cells.Item(17, 8).Value = [double]"5.88727696920018"

$ws.Cells.Item(18, 1).Value = [double]"0.783853560058326"
$ws.Cells.Item(18, 2).Value = [double]"2.775251653342024"
$ws.Cells.Item(18, 3).Value = [double]"0.7032173951324423"
$ws.Cells.Item(18, 4).Value = [double]"-0.4754455809308937"
$ws.Cells.Item(18, 5).Value = [double]"-27.24102517548582"
$ws.Cells.Item(18, 6).Value = [double]"1.500000014390762"
$ws.Cells.Item(18, 7).Value = [double]"23.15818705638852"
$ws.Cells.Item(18, 8).Value = [double]"5.887276969720969"

$ws.Cells.Item(19, 1).Value = [double]"0.8328444075619713"
$ws.Cells.Item(19, 2).Value = [double]"2.958997190743494"
$ws.Cells.Item(19, 3).Value = [double]"0.6086186385217364"
$ws.Cells.Item(19, 4).Value = [double]"-0.4267147203047844"
$ws.Cells.Item(19, 5).Value = [double]"-24.44895252956952"
$ws.Cells.Item(19, 6).Value = [double]"1.500000014507612"
$ws.Cells.Item(19, 7).Value = [double]"27.3686402918748"
$ws.Cells.Item(19, 8).Value = [double]"5.887276970103439"

$ws.Cells.Item(20, 1).Value = [double]"0.8818352550656167"
$ws.Cells.Item(20, 2).Value = [double]"3.147132655877436"
$ws.Cells.Item(20, 3).Value = [double]"0.5230827160909686"
$ws.Cells.Item(20, 4).Value = [double]"-0.3677395764706727"
$ws.Cells.Item(20, 5).Value = [double]"-21.06992569169794"
$ws.Cells.Item(20, 6).Value = [double]"1.50000001459516"
$ws.Cells.Item(20, 7).Value = [double]"31.57909352768907"
$ws.Cells.Item(20, 8).Value = [double]"5.887276970396274"

$ws.Cells.Item(21, 1).Value = [double]"0.9308261025692621"
$ws.Cells.Item(21, 2).Value = [double]"3.33998261155157"
$ws.Cells.Item(21, 3).Value = [double]"0.4487843859365749"
$ws.Cells.Item(21, 4).Value = [double]"-0.2977084102017189"
$ws.Cells.Item(21, 5).Value = [double]"-17.05743543010795"
$ws.Cells.Item(21, 6).Value = [double]"1.50000001466317"
$ws.Cells.Item(21, 7).Value = [double]"35.78954676374908"
$ws.Cells.Item(21, 8).Value = [double]"1.989055301207893"

$ws.Cells.Item(22, 1).Value = [double]"0.9798169500729075"
$ws.Cells.Item(22, 2).Value = [double]"3.537558803945962"
$ws.Cells.Item(22, 3).Value = [double]"0.3881626450659849"
$ws.Cells.Item(22, 4).Value = [double]"-0.2155731427332745"
$ws.Cells.Item(22, 5).Value = [double]"-12.35143125498792"
$ws.Cells.Item(22, 6).Value = [double]"-2.008874664160104e-24"
$ws.Cells.Item(22, 7).Value = [double]"40"
$ws.Cells.Item(22, 8).Value = [double]"2.163573273794861e-16"

Write-Host "Applied all updates"
